$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking columns (D, E, G) so values
# such as "283.20", "1.85%", "12" are retained as text, matching the source data.
$ws.Range("D2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "283.20"
$ws.Range("E2").Value = "1.85%"
$ws.Range("G2").Value = "12"
$ws.Range("D3").Value = "28.29"
$ws.Range("E3").Value = "4.00%"
$ws.Range("G3").Value = "12"
$ws.Range("D4").Value = "5.018"
$ws.Range("E4").Value = "3.09%"
$ws.Range("G4").Value = "12"
$ws.Range("D5").Value = "0.06520"
$ws.Range("E5").Value = "1.64%"
$ws.Range("G5").Value = "12"
$ws.Range("D6").Value = "7.244"
$ws.Range("E6").Value = "3.25%"
$ws.Range("G6").Value = "12"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "3.360"
$ws.Range("E7").Value = "1.67%"
$ws.Range("G7").Value = "12"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "1.409"
$ws.Range("E8").Value = "17.05%"
$ws.Range("G8").Value = "12"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9176"
$ws.Range("E9").Value = "3.38%"
$ws.Range("G9").Value = "12"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1549"
$ws.Range("E10").Value = "0.42%"
$ws.Range("G10").Value = "12"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "0.06549"
$ws.Range("E11").Value = "26.82%"
$ws.Range("G11").Value = "12"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07541"
$ws.Range("E12").Value = "0.54%"
$ws.Range("G12").Value = "12"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.02756"
$ws.Range("E13").Value = "-4.61%"
$ws.Range("G13").Value = "12"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.08959"
$ws.Range("E14").Value = "-0.16%"
$ws.Range("G14").Value = "12"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001587"
$ws.Range("E15").Value = "1.25%"
$ws.Range("G15").Value = "12"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "0.0006377"
$ws.Range("E16").Value = "-0.11%"
$ws.Range("G16").Value = "12"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.006050"
$ws.Range("E17").Value = "-1.82%"
$ws.Range("G17").Value = "12"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.439"
$ws.Range("E18").Value = "-1.06%"
$ws.Range("G18").Value = "12"
$ws.Range("D19").Value = "2.242"
$ws.Range("E19").Value = "-0.23%"
$ws.Range("G19").Value = "12"
$ws.Range("D20").Value = "0.3186"
$ws.Range("E20").Value = "3.01%"
$ws.Range("G20").Value = "12"
$ws.Range("D21").Value = "0.1279"
$ws.Range("E21").Value = "-4.69%"
$ws.Range("G21").Value = "12"
$ws.Range("D22").Value = "3.974"
$ws.Range("E22").Value = "1.48%"
$ws.Range("G22").Value = "12"
$ws.Range("D23").Value = "0.1542"
$ws.Range("E23").Value = "1.58%"
$ws.Range("G23").Value = "12"
$ws.Range("D24").Value = "0.04419"
$ws.Range("E24").Value = "0.24%"
$ws.Range("G24").Value = "12"
$ws.Range("D25").Value = "0.001182"
$ws.Range("E25").Value = "0.39%"
$ws.Range("G25").Value = "12"
$ws.Range("D26").Value = "0.004438"
$ws.Range("E26").Value = "14.11%"
$ws.Range("G26").Value = "12"
$ws.Range("E27").Value = "1.56%"
$ws.Range("G27").Value = "12"
$ws.Range("D28").Value = "0.0001617"
$ws.Range("E28").Value = "-1.58%"
$ws.Range("G28").Value = "12"
$ws.Range("G29").Value = "12"
$ws.Range("G30").Value = "12"
$ws.Range("G31").Value = "12"
$ws.Range("G32").Value = "12"
$ws.Range("G33").Value = "12"
$ws.Range("G34").Value = "12"
$ws.Range("G35").Value = "12"
$ws.Range("G36").Value = "12"
$ws.Range("G37").Value = "12"
$ws.Range("G38").Value = "12"
$ws.Range("G39").Value = "12"
$ws.Range("D40").Value = "0.04126"
$ws.Range("E40").Value = "0.24%"
$ws.Range("G40").Value = "12"
$ws.Range("D41").Value = "0.006667"
$ws.Range("E41").Value = "-2.00%"
$ws.Range("G41").Value = "12"
$ws.Range("D42").Value = "0.1230"
$ws.Range("E42").Value = "4.74%"
$ws.Range("G42").Value = "12"
$ws.Range("D43").Value = "0.002178"
$ws.Range("E43").Value = "15.24%"
$ws.Range("G43").Value = "12"
$ws.Range("D44").Value = "0.01207"
$ws.Range("E44").Value = "3.81%"
$ws.Range("G44").Value = "12"
$ws.Range("D45").Value = "0.00005648"
$ws.Range("E45").Value = "5.60%"
$ws.Range("G45").Value = "12"
$ws.Range("E46").Value = "20.74%"
$ws.Range("G46").Value = "12"
$ws.Range("E47").Value = "0.02%"
$ws.Range("G47").Value = "12"
$ws.Range("G48").Value = "12"
$ws.Range("G49").Value = "12"
$ws.Range("G50").Value = "12"
$ws.Range("G51").Value = "12"

# Restore default (General) style on the numeric-looking columns so no
# residual text-number-format style is left behind on the cells.
$ws.Range("D2:G51").Style = "Normal"
